$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(-17.48052200551075, 1.715765534478213, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075),
  @(-17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075),
  @(-17.48052200551075, 2.148181613717205, 2.013093073539469, -17.48052200551075, 3.678283039360347, -17.48052200551075, 1.797829893983427, -17.48052200551075, 1.406451625362893, -17.48052200551075),
  @(-17.48052200551075, 1.860868682517347, -17.48052200551075, -17.48052200551075, -17.48052200551075, 2.38528525616807, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075),
  @(-17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075),
  @(2.664973606895619, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075),
  @(-17.48052200551075, -17.48052200551075, -17.48052200551075, 2.125978997129963, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075),
  @(3.771647345949205, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075),
  @(-17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, 2.116259214142572),
  @(-17.48052200551075, -17.48052200551075, -17.48052200551075, 2.802472494383756, -17.48052200551075, 3.48032345304084, -17.48052200551075, -17.48052200551075, -17.48052200551075, 2.551837220762657),
  @(-17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075),
  @(-17.48052200551075, -17.48052200551075, -17.48052200551075, 2.500845984578258, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, 0.6174329468181577, 1.664676290586494),
  @(-17.48052200551075, -17.48052200551075, 1.071430278967461, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, 2.124713333925178),
  @(-17.48052200551075, -17.48052200551075, 0.9927437729820173, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075),
  @(-17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, 1.360861134042961, -17.48052200551075),
  @(-17.48052200551075, 2.377472653277211, 2.810973610821283, -17.48052200551075, -17.48052200551075, -17.48052200551075, 1.923432433216881, 4.321920600471146, 3.079793100940034, -17.48052200551075),
  @(-17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, -17.48052200551075, 2.621998398615541, -17.48052200551075, 2.26051204685694, -17.48052200551075),
  @(-17.48052200551075, -17.48052200551075, 1.399052270304218, -17.48052200551075, -17.48052200551075, -17.48052200551075, 0.663782728949743, -17.48052200551075, -17.48052200551075, -17.48052200551075),
  @(-17.48052200551075, 0.7494280430386079, 1.14868162180884, -17.48052200551075, 2.847601753106291, -17.48052200551075, 1.241087895934416, -17.48052200551075, -17.48052200551075, 1.181955145728686),
  @(-17.48052200551075, 0.8260960070815589, -17.48052200551075, 1.584191818331799, -17.48052200551075, 1.853989165875281, 1.392696788470908, -17.48052200551075, -17.48052200551075, -17.48052200551075)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $row[$c]
    }
}
